$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of bulk data. Values that look numeric ("8", "34545") must be
# stored as literal text — the apostrophe prefix forces text type even for
# digit-only strings (matching the source data's inlineStr "STP ID" /
# "Batch Number" cells), and resetting the style back to "Normal" clears
# the transient quote-prefix formatting so no stray style index is left
# behind on the new cells.
$ws.Range("A9").Value = "form_8.pdf"

$ws.Range("B9").Value = "'8"
$ws.Range("B9").Style = "Normal"

$ws.Range("C9").Value = "Paracetamole"

$ws.Range("D9").Value = "'34545"
$ws.Range("D9").Style = "Normal"

# Remaining mid-row columns are present but empty (empty text, not blank).
$ws.Range("E9:O9").Value = "'"
$ws.Range("E9:O9").Style = "Normal"

$ws.Range("P9").Value = "TEST"
